# [23rd June 23]:- login Script designed
$wb = $excel.ActiveWorkbook

# Sheet2 ("Sheet2") loses the active/selected tab flag and its selection
# moves from A4 to B5.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B5").Select() | Out-Null

# Sheet3 ("abc") is renamed to "DDF", gets a new header cell, becomes the
# active sheet, is zoomed to 205% and the selection moves to A3.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "DDF"
$ws3.Range("A1").Value = "Jan batch"

$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 205
$ws3.Range("A3").Select() | Out-Null
